# Automatische test-sync: 2025-07-31 22:04:50
# Appends a new test-mail log row to the "Logs" sheet and refreshes the
# "Dashboard" summary sheet to match (swap the "Productinformatie" /
# "Intern verzoek / Actie voor medewerker" category counts).

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 21 ------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(21, 1).Value  = "Bel jij klant Jansen even?"
$logs.Cells.Item(21, 2).Value  = "mailmind.test@zohomail.eu"
$logs.Cells.Item(21, 3).Value  = "Testmail #19: Bel jij klant Jansen even?"
$logs.Cells.Item(21, 4).Value  = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item(21, 5).Value  = "Dank voor je bericht. We pakken dit intern op en houden je op de hoogte."
$logs.Cells.Item(21, 6).Value  = "2025-07-31 22:04:06"
$logs.Cells.Item(21, 7).Value  = "Ja"
$logs.Cells.Item(21, 8).Value  = "Ja"
$logs.Cells.Item(21, 9).Value  = "Nee"
$logs.Cells.Item(21, 10).Value = "Nee"

# --- Dashboard sheet: update category summary rows ---------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(3, 1).Value = "Intern verzoek / Actie voor medewerker"
$dash.Cells.Item(3, 2).Value = 5
$dash.Cells.Item(4, 1).Value = "Productinformatie"
$dash.Cells.Item(4, 2).Value = 5

# --- Logs sheet: extend conditional formatting ranges to include row 21 ------
$logs.Range("D2:D20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D21"))
$logs.Range("G2:G20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G21"))
$logs.Range("H2:H20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H21"))
$logs.Range("I2:I20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I21"))
$logs.Range("J2:J20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J21"))
